# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 to match the latest snapshot.
# Column D cells are forced to Text format before assignment so that
# numeric-looking price strings (e.g. "0.9994") remain text, matching
# the original inlineStr cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.056.73"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.754.88"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.13"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3776"
$ws.Range("E7").Value = "  -4.19%  "
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.37"
$ws.Range("E9").Value = "  -5.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07196"
$ws.Range("E11").Value = "  -4.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.54"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.135"
$ws.Range("E14").Value = "  -5.94%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.757.20"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06583"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.48"
$ws.Range("E19").Value = "  -5.78%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.87"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.235"
$ws.Range("E22").Value = "  -5.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.081.63"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  -5.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.64"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.79"
$ws.Range("E27").Value = "  -7.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.319"
$ws.Range("E28").Value = "  -8.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.956.25"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.65"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.246"
$ws.Range("E31").Value = "  -15.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.019"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.782"
$ws.Range("E33").Value = "  -7.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08717"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.20"
$ws.Range("E35").Value = "  -7.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6680"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02325"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06181"
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.155"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.215"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.446"
$ws.Range("E42").Value = "  -10.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.021"
$ws.Range("E43").Value = "  -6.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.59"
$ws.Range("E45").Value = "  -6.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.834"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6046"
$ws.Range("E47").Value = "  -6.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.67"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("E49").Value = "  -6.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07152"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.169"
$ws.Range("E51").Value = "  +0.53%  "